$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '37.058.30'
Set-TextValue "E2" '  -1.81%  '
Set-TextValue "D3" '2.017.09'
Set-TextValue "E3" '  -2.98%  '
Set-TextValue "E4" '  +0.02%  '
Set-TextValue "D5" '226.12'
Set-TextValue "E5" '  -2.77%  '
Set-TextValue "D6" '0.604'
Set-TextValue "E6" '  -3.05%  '
Set-TextValue "E7" '  +0.00%  '
Set-TextValue "D8" '54.72'
Set-TextValue "E8" '  -5.86%  '
Set-TextValue "D9" '0.376'
Set-TextValue "E9" '  -4.08%  '
Set-TextValue "D10" '0.0784'
Set-TextValue "E10" '  +0.38%  '
Set-TextValue "D11" '0.103'
Set-TextValue "E11" '  -5.35%  '
Set-TextValue "D12" '2.312.17'
Set-TextValue "E12" '  -3.05%  '
Set-TextValue "D13" '14.11'
Set-TextValue "E13" '  -5.14%  '
Set-TextValue "D14" '20.18'
Set-TextValue "E14" '  -4.93%  '
Set-TextValue "E15" '  -3.64%  '
Set-TextValue "D16" '5.13'
Set-TextValue "E16" '  -3.82%  '
Set-TextValue "D17" '2.019.96'
Set-TextValue "E17" '  -2.69%  '
Set-TextValue "D18" '36.994.97'
Set-TextValue "E18" '  -1.81%  '
Set-TextValue "D19" '6.21'
Set-TextValue "E19" '  +0.63%  '
Set-TextValue "D20" '68.89'
Set-TextValue "E20" '  -1.90%  '
Set-TextValue "D21" '0.0₃0816'
Set-TextValue "E21" '  -1.85%  '
Set-TextValue "D22" '223.03'
Set-TextValue "E22" '  -2.12%  '
Set-TextValue "E23" '  -0.01%  '
Set-TextValue "E24" '  +1.83%  '
Set-TextValue "E25" '  -7.67%  '
Set-TextValue "D26" '166.20'
Set-TextValue "E26" '  -1.96%  '
Set-TextValue "D27" '9.15'
Set-TextValue "E27" '  -7.89%  '
Set-TextValue "B28" 'ImmutableX'
Set-TextValue "C28" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D28" '1.35'
Set-TextValue "E28" '  -2.38%  '
Set-TextValue "B29" 'EthereumClassic'
Set-TextValue "C29" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D29" '18.68'
Set-TextValue "E29" '  -3.58%  '
Set-TextValue "E30" '  -6.63%  '
Set-TextValue "E31" '  -3.71%  '
Set-TextValue "E32" '  -2.29%  '
Set-TextValue "D33" '0.0611'
Set-TextValue "E33" '  -3.05%  '
Set-TextValue "D34" '4.41'
Set-TextValue "E34" '  -5.07%  '
Set-TextValue "D35" '2.34'
Set-TextValue "E35" '  -7.41%  '
Set-TextValue "E36" '  +0.72%  '
Set-TextValue "E37" '  -0.08%  '
Set-TextValue "D38" '3.15'
Set-TextValue "E38" '  -5.46%  '
Set-TextValue "D39" '5.29'
Set-TextValue "E39" '  -1.19%  '
Set-TextValue "D40" '1.476.55'
Set-TextValue "E40" '  -0.74%  '
Set-TextValue "E41" '  -5.37%  '
Set-TextValue "D42" '94.82'
Set-TextValue "E42" '  -3.67%  '
Set-TextValue "D43" '0.0912'
Set-TextValue "E43" '  -4.86%  '
Set-TextValue "B44" 'InjectiveProtocol'
Set-TextValue "C44" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D44" '16.23'
Set-TextValue "E44" '  -3.99%  '
Set-TextValue "B45" 'HuobiToken'
Set-TextValue "C45" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D45" '2.77'
Set-TextValue "E45" '  -5.22%  '
Set-TextValue "E46" '  -5.99%  '
Set-TextValue "D47" '7.19'
Set-TextValue "E47" '  -1.28%  '
Set-TextValue "E48" '  -3.37%  '
Set-TextValue "E49" '  -1.75%  '
Set-TextValue "D50" '2.199.68'
Set-TextValue "E50" '  -3.07%  '
Set-TextValue "D51" '44.16'
Set-TextValue "E51" '  -4.07%  '
